$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F header: copy the header formatting from E1 (bold/border/wrap
# header style used by the whole header row), then set the new label text.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Расширенный комментарий"

# New column F width (~24 "characters" once Excel's internal padding is
# applied back when the workbook is re-saved).
$ws.Columns("F").ColumnWidth = 23.14

# F5 gets the same "empty but bordered" look as E5 (vertical-top/wrap style).
$ws.Range("E5").Copy()
$ws.Range("F5").PasteSpecial(-4122)

# Match the workbook's recorded selection/active cell after the edit.
$ws.Range("F1").Select()
